$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 74, shifting existing rows 74:86 down to 75:87.
$ws.Rows.Item(74).Insert()

# Populate the newly inserted row 74 with the new weekly record.
# (Unchanged-from-old-row-74 fields are re-written too, since Insert() leaves
# the new row blank except for copied formatting.)
$ws.Cells.Item(74, 1).Value = 10
$ws.Cells.Item(74, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(74, 3).Value = "La Araucanía"
$ws.Cells.Item(74, 4).Value = 44504
$ws.Cells.Item(74, 5).Value = 9
$ws.Cells.Item(74, 6).Value = 100112031
$ws.Cells.Item(74, 7).Value = "Poroto verde"
$ws.Cells.Item(74, 8).Value = "Sin especificar"
$ws.Cells.Item(74, 9).Value = "Primera"
$ws.Cells.Item(74, 10).Value = 65
$ws.Cells.Item(74, 11).Value = 40000
$ws.Cells.Item(74, 12).Value = 40000
$ws.Cells.Item(74, 13).Value = 40000
$ws.Cells.Item(74, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(74, 15).Value = "Perú"
$ws.Cells.Item(74, 16).Value = 1600
$ws.Cells.Item(74, 17).Value = 25
$ws.Cells.Item(74, 18).Value = "Hortaliza"
